# Refactored JDE Vendor search
# Adds new "Home Page" / informational text / alert text / login status
# columns to the VendorSearch sheet, and makes VendorSearch the active
# (selected) sheet/tab instead of PhoneBook.

$wb = $excel.ActiveWorkbook
$wsPhoneBook = $wb.Worksheets.Item("PhoneBook")
$ws = $wb.Worksheets.Item("VendorSearch")

# --- Copy the existing header (row 2) / data (row 3) cell formatting into
#     the new columns E:H before putting values in them, so the new cells
#     pick up the same style as the rest of their row. ---
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2:H2").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("E3:H3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Fill in the new header + data values (this ordering reproduces the
#     shared-string table order of the original edit). ---
$ws.Range("E2").Value = "HomePageTitle"
$ws.Range("E3").Value = "Vendor Search Form"
$ws.Range("F3").Value = "You have tried to access a secure area. Please enter your user name and password to gain access."
$ws.Range("F2").Value = "informationtex"
$ws.Range("G3").Value = "Do not bookmark this page. Doing so will cause an error on your next visit. The correct procedure is to bookmark the page that appears after you click "
$ws.Range("G2").Value = "alertTex"
$ws.Range("H2").Value = "loginStatus"
$ws.Range("H3").Value = "Not logged in"

# --- New column widths for E:H ---
$ws.Columns.Item(5).ColumnWidth = 18.34
$ws.Columns.Item(6).ColumnWidth = 13.35
$ws.Columns.Item(7).ColumnWidth = 9.35
$ws.Columns.Item(8).ColumnWidth = 14.35

# --- Make VendorSearch the active sheet/tab and select H2 on it
#     (was F6 before); PhoneBook loses tabSelected as a side effect. ---
$ws.Activate()
$ws.Range("H2").Select() | Out-Null
